# Update cryptos list (prices, volume %, and a few reordered rows)
# per refreshed data from coinranking.com feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.882.88'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.843.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '704.26'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.97'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.841.97'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.34'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.50'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.492.87'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.835.97'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.941.93'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.33'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.67'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.30'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.18'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.06%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.182'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.39'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.799.86'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.14'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.36'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.02'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.50%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.79%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.19%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("B45").Value = 'FLOKI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000311'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.71%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '163.20'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.80'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '413.22'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.25%  '
$ws.Range("B49").Value = 'TheGraph'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.299'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.62'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.18'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.02%  '
